$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planning")

# --- Rows 2/3/4 ("Folded Sheet  4p") -----------------------------------
# Cyclically rotate the Activity/ScheduledResource/PlannedQty/
# MISWorkCenter/PlannedResource cells among rows 2, 3 and 4:
#   new row2 <- old row4, new row3 <- old row2, new row4 <- old row3
$rotCols = @("B", "D", "G", "O", "P")

$row2 = @{}
$row3 = @{}
$row4 = @{}
foreach ($col in $rotCols) {
    $row2[$col] = $ws.Range($col + "2").Value()
    $row3[$col] = $ws.Range($col + "3").Value()
    $row4[$col] = $ws.Range($col + "4").Value()
}

foreach ($col in $rotCols) {
    $ws.Range($col + "2").Value = $row4[$col]
    $ws.Range($col + "3").Value = $row2[$col]
    $ws.Range($col + "4").Value = $row3[$col]
}

# --- Rows 22/23 ("Plate - Folded Sheet 1  4p") --------------------------
# Swap the Activity/ScheduledResource/MISWorkCenter/PlannedResource cells
# between the two rows so row 22 becomes the "Plate burn / Plate Making"
# step and row 23 becomes the "OKTP / Ok to Plate" step - matching the
# pattern used by the other triplets (rows 20/21, 24/25).
$swapCols = @("B", "D", "O", "P")

foreach ($col in $swapCols) {
    $cell22 = $ws.Range($col + "22")
    $cell23 = $ws.Range($col + "23")

    $val22 = $cell22.Value()
    $val23 = $cell23.Value()

    $cell22.Value = $val23
    $cell23.Value = $val22
}
